$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.503.61'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').Value = '1.791.53'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').Value = '306.02'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').Value = '0.4276'
$ws.Range('E7').Value = '  +2.02%  '
$ws.Range('D8').Value = '0.3619'
$ws.Range('E8').Value = '  +1.12%  '
$ws.Range('D9').Value = '0.07143'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('D10').Value = '0.8549'
$ws.Range('E10').Value = '  +0.92%  '
$ws.Range('D11').Value = '20.60'
$ws.Range('E11').Value = '  +2.05%  '
$ws.Range('D12').Value = '1.870.23'
$ws.Range('E12').Value = '  +4.17%  '
$ws.Range('D13').Value = '6.517'
$ws.Range('E13').Value = '  +2.53%  '
$ws.Range('D14').Value = '5.282'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = '0.06876'
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '79.76'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '0.000008832'
$ws.Range('E18').Value = '  +1.50%  '
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').Value = '15.00'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').Value = '26.531.67'
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('D22').Value = '5.140'
$ws.Range('E22').Value = '  +1.67%  '
$ws.Range('D23').Value = '11.04'
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('D24').Value = '2.080.76'
$ws.Range('E24').Value = '  +3.45%  '
$ws.Range('D25').Value = '152.06'
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('D26').Value = '1.825'
$ws.Range('E26').Value = '  -5.65%  '
$ws.Range('D27').Value = '18.12'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = '5.165'
$ws.Range('E28').Value = '  +2.92%  '
$ws.Range('D29').Value = '1.893'
$ws.Range('E29').Value = '  +14.98%  '
$ws.Range('D30').Value = '114.88'
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('D31').Value = '0.08900'
$ws.Range('E31').Value = '  -1.17%  '
$ws.Range('D32').Value = '0.7451'
$ws.Range('E32').Value = '  +3.28%  '
$ws.Range('E33').Value = '  +5.53%  '
$ws.Range('D34').Value = '4.351'
$ws.Range('E34').Value = '  +1.08%  '
$ws.Range('D35').Value = '2.747'
$ws.Range('E35').Value = '  -3.86%  '
$ws.Range('D36').Value = '1.002'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').Value = '1.112'
$ws.Range('E37').Value = '  +3.33%  '
$ws.Range('D38').Value = '0.05158'
$ws.Range('E38').Value = '  +0.70%  '
$ws.Range('D39').Value = '0.01900'
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('D40').Value = '0.4991'
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('D41').Value = '0.1622'
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').Value = '2.611'
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('D43').Value = '6.422'
$ws.Range('E43').Value = '  +7.44%  '
$ws.Range('D44').Value = '8.215'
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('D45').Value = '105.72'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('D46').Value = '10.28'
$ws.Range('E46').Value = '  +0.80%  '
$ws.Range('D47').Value = '1.002'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '1.642'
$ws.Range('E48').Value = '  +2.37%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.06207'
$ws.Range('E49').Value = '  -1.50%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').Value = '0.4501'
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('D51').Value = '1.776'
$ws.Range('E51').Value = '  +4.65%  '
